$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "37-29="
$t.Cell(1, 2).Range.Text = "28+69="
$t.Cell(1, 3).Range.Text = "67+28="
$t.Cell(1, 4).Range.Text = "42-34="
$t.Cell(1, 5).Range.Text = "54-29="
$t.Cell(2, 1).Range.Text = "52-33="
$t.Cell(2, 2).Range.Text = "45+16="
$t.Cell(2, 3).Range.Text = "91-8="
$t.Cell(2, 4).Range.Text = "61-53="
$t.Cell(2, 5).Range.Text = "46+28="
$t.Cell(3, 1).Range.Text = "23-5="
$t.Cell(3, 2).Range.Text = "68+4="
$t.Cell(3, 3).Range.Text = "54-29="
$t.Cell(3, 4).Range.Text = "74-26="
$t.Cell(3, 5).Range.Text = "63-18="
$t.Cell(4, 1).Range.Text = "45+18="
$t.Cell(4, 2).Range.Text = "25+58="
$t.Cell(4, 3).Range.Text = "37+37="
$t.Cell(4, 4).Range.Text = "57+19="
$t.Cell(4, 5).Range.Text = "70-12="
$t.Cell(5, 1).Range.Text = "23-8="
$t.Cell(5, 2).Range.Text = "50-5="
$t.Cell(5, 3).Range.Text = "90-35="
$t.Cell(5, 4).Range.Text = "83-55="
$t.Cell(5, 5).Range.Text = "69+13="
$t.Cell(6, 1).Range.Text = "81-37="
$t.Cell(6, 2).Range.Text = "82-59="
$t.Cell(6, 3).Range.Text = "23-5="
$t.Cell(6, 4).Range.Text = "17+48="
$t.Cell(6, 5).Range.Text = "12-6="
$t.Cell(7, 1).Range.Text = "63+28="
$t.Cell(7, 2).Range.Text = "9+45="
$t.Cell(7, 3).Range.Text = "61-58="
$t.Cell(7, 4).Range.Text = "85-49="
$t.Cell(7, 5).Range.Text = "37+36="
$t.Cell(8, 1).Range.Text = "8+87="
$t.Cell(8, 2).Range.Text = "2+39="
$t.Cell(8, 3).Range.Text = "37+27="
$t.Cell(8, 4).Range.Text = "30-15="
$t.Cell(8, 5).Range.Text = "18+69="
$t.Cell(9, 1).Range.Text = "58+35="
$t.Cell(9, 2).Range.Text = "84-77="
$t.Cell(9, 3).Range.Text = "18+26="
$t.Cell(9, 4).Range.Text = "31-25="
$t.Cell(9, 5).Range.Text = "40-13="
$t.Cell(10, 1).Range.Text = "56+38="
$t.Cell(10, 2).Range.Text = "27+48="
$t.Cell(10, 3).Range.Text = "8+3="
$t.Cell(10, 4).Range.Text = "39+26="
$t.Cell(10, 5).Range.Text = "65+19="
$t.Cell(11, 1).Range.Text = "6+18="
$t.Cell(11, 2).Range.Text = "61-13="
$t.Cell(11, 3).Range.Text = "40-25="
$t.Cell(11, 4).Range.Text = "8+49="
$t.Cell(11, 5).Range.Text = "7+88="
$t.Cell(12, 1).Range.Text = "70-47="
$t.Cell(12, 2).Range.Text = "92-14="
$t.Cell(12, 3).Range.Text = "51-29="
$t.Cell(12, 4).Range.Text = "46-27="
$t.Cell(12, 5).Range.Text = "61-28="
$t.Cell(13, 1).Range.Text = "86-69="
$t.Cell(13, 2).Range.Text = "60-53="
$t.Cell(13, 3).Range.Text = "43+18="
$t.Cell(13, 4).Range.Text = "16+55="
$t.Cell(13, 5).Range.Text = "49+5="
$t.Cell(14, 1).Range.Text = "36-28="
$t.Cell(14, 2).Range.Text = "2+49="
$t.Cell(14, 3).Range.Text = "81-77="
$t.Cell(14, 4).Range.Text = "59+14="
$t.Cell(14, 5).Range.Text = "33+48="
$t.Cell(15, 1).Range.Text = "29+47="
$t.Cell(15, 2).Range.Text = "47+5="
$t.Cell(15, 3).Range.Text = "18+54="
$t.Cell(15, 4).Range.Text = "95-89="
$t.Cell(15, 5).Range.Text = "81-23="
$t.Cell(16, 1).Range.Text = "34-7="
$t.Cell(16, 2).Range.Text = "55-26="
$t.Cell(16, 3).Range.Text = "19+72="
$t.Cell(16, 4).Range.Text = "60-12="
$t.Cell(16, 5).Range.Text = "68+19="
$t.Cell(17, 1).Range.Text = "54-28="
$t.Cell(17, 2).Range.Text = "80-56="
$t.Cell(17, 3).Range.Text = "8+26="
$t.Cell(17, 4).Range.Text = "29+4="
$t.Cell(17, 5).Range.Text = "28+6="
$t.Cell(18, 1).Range.Text = "35+18="
$t.Cell(18, 2).Range.Text = "90-86="
$t.Cell(18, 3).Range.Text = "57+29="
$t.Cell(18, 4).Range.Text = "41-7="
$t.Cell(18, 5).Range.Text = "67+9="
$t.Cell(19, 1).Range.Text = "83-4="
$t.Cell(19, 2).Range.Text = "39+15="
$t.Cell(19, 3).Range.Text = "36+57="
$t.Cell(19, 4).Range.Text = "5+49="
$t.Cell(19, 5).Range.Text = "26+35="
$t.Cell(20, 1).Range.Text = "59+24="
$t.Cell(20, 2).Range.Text = "70-34="
$t.Cell(20, 3).Range.Text = "37+27="
$t.Cell(20, 4).Range.Text = "67-18="
$t.Cell(20, 5).Range.Text = "28+46="
